$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the Comercializadora del Agro de Limarí
# Chirimoya block (row 53), pushing the existing rows 53-70 down to 56-73.
$ws.Rows("53:55").Insert()

# New weekly entries for fecha 2021-11-10 (serial 44510): Especial, Primera,
# Segunda (no "Extra (doble especial)" entry for this date).

# Row 53 - Especial
$ws.Range("A53").Value = 2
$ws.Range("B53").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44510
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107002
$ws.Range("J53").Value = "Chirimoya"
$ws.Range("K53").Value = "Cultivar IV Región"
$ws.Range("L53").Value = "Especial"
$ws.Range("M53").Value = 360
$ws.Range("N53").Value = 1600
$ws.Range("O53").Value = 1700
$ws.Range("P53").Value = 1650
$ws.Range("Q53").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R53").Value = "Provincia de Limarí"
$ws.Range("S53").Value = 1650
$ws.Range("T53").Value = 1

# Row 54 - Primera
$ws.Range("A54").Value = 2
$ws.Range("B54").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44510
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100107
$ws.Range("H54").Value = "Otros"
$ws.Range("I54").Value = 100107002
$ws.Range("J54").Value = "Chirimoya"
$ws.Range("K54").Value = "Cultivar IV Región"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 360
$ws.Range("N54").Value = 1200
$ws.Range("O54").Value = 1300
$ws.Range("P54").Value = 1250
$ws.Range("Q54").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R54").Value = "Provincia de Limarí"
$ws.Range("S54").Value = 1250
$ws.Range("T54").Value = 1

# Row 55 - Segunda
$ws.Range("A55").Value = 2
$ws.Range("B55").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44510
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100107
$ws.Range("H55").Value = "Otros"
$ws.Range("I55").Value = 100107002
$ws.Range("J55").Value = "Chirimoya"
$ws.Range("K55").Value = "Cultivar IV Región"
$ws.Range("L55").Value = "Segunda"
$ws.Range("M55").Value = 360
$ws.Range("N55").Value = 1000
$ws.Range("O55").Value = 1100
$ws.Range("P55").Value = 1050
$ws.Range("Q55").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R55").Value = "Provincia de Limarí"
$ws.Range("S55").Value = 1050
$ws.Range("T55").Value = 1
